# Updated cryptos list: refresh Price (D) / Volume(1h) (E) figures, and
# restore the correct Coin/Link/Price/Volume rows for entries #32 and #33
# which had been swapped (BinanceUSD <-> ImmutableX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell reference -> new text value, taken verbatim from the refreshed feed
$updates = [ordered]@{
    'D2' = '35.331.09'
    'E2' = '  +0.39%  '
    'D3' = '1.909.43'
    'E3' = '  +2.66%  '
    'E4' = '  -0.36%  '
    'D5' = '245.46'
    'E5' = '  +2.46%  '
    'E6' = '  +6.10%  '
    'E7' = '  -0.39%  '
    'D8' = '41.17'
    'E8' = '  -2.66%  '
    'D9' = '0.349'
    'E9' = '  +5.79%  '
    'D10' = '52.70'
    'E10' = '  +12.41%  '
    'D11' = '0.0716'
    'E11' = '  +3.22%  '
    'D12' = '0.0992'
    'E12' = '  +0.48%  '
    'D13' = '2.189.22'
    'E13' = '  +2.86%  '
    'D14' = '12.09'
    'E14' = '  +5.19%  '
    'D15' = '0.701'
    'E15' = '  +3.40%  '
    'D16' = '1.914.97'
    'E16' = '  +2.97%  '
    'E17' = '  +2.86%  '
    'D18' = '35.334.67'
    'E18' = '  +0.50%  '
    'D19' = '72.28'
    'E19' = '  +3.40%  '
    'D20' = '0.0₃0828'
    'E20' = '  +3.83%  '
    'D21' = '239.52'
    'E21' = '  -0.45%  '
    'D22' = '12.51'
    'E22' = '  +2.21%  '
    'D23' = '4.84'
    'E23' = '  +2.06%  '
    'E24' = '  -0.41%  '
    'D25' = '2.29'
    'E25' = '  +1.25%  '
    'E26' = '  +22.87%  '
    'D27' = '169.55'
    'E27' = '  +0.51%  '
    'D28' = '8.50'
    'E28' = '  +6.20%  '
    'D29' = '18.44'
    'E29' = '  +4.46%  '
    'D30' = '0.126'
    'E30' = '  +1.86%  '
    'E31' = '  +3.70%  '
    'D32' = '0.0567'
    'E32' = '  +1.70%  '
    'B33' = 'BinanceUSD'
    'C33' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D33' = '1.02'
    'E33' = '  +0.45%  '
    'B34' = 'ImmutableX'
    'C34' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D34' = '0.938'
    'E34' = '  +14.99%  '
    'D35' = '4.12'
    'E35' = '  +2.49%  '
    'D36' = '1.74'
    'E36' = '  -4.08%  '
    'E37' = '  +1.24%  '
    'E38' = '  +0.83%  '
    'D39' = '1.11'
    'E39' = '  +1.46%  '
    'D40' = '0.0653'
    'E40' = '  +10.12%  '
    'D41' = '0.0208'
    'E41' = '  +3.88%  '
    'D42' = '16.31'
    'E42' = '  +8.96%  '
    'D43' = '89.90'
    'E43' = '  -0.20%  '
    'D44' = '1.338.15'
    'E44' = '  -0.42%  '
    'E45' = '  +3.10%  '
    'D46' = '48.45'
    'E46' = '  +39.23%  '
    'E47' = '  +1.98%  '
    'E48' = '  -0.31%  '
    'D49' = '6.58'
    'E49' = '  -0.31%  '
    'D50' = '2.096.30'
    'E50' = '  +2.70%  '
    'E51' = '  +3.67%  '
}

foreach ($cellRef in $updates.Keys) {
    $text = $updates[$cellRef]
    $range = $ws.Range($cellRef)

    # Column D holds price strings such as "245.46" or "52.70". Handed a
    # plain-looking decimal, Excel auto-converts the cell to a Number
    # (dropping the trailing zero / exact text), so force text storage the
    # same way a user typing `245.46` into the grid would: a leading
    # apostrophe. Values that are not plain decimals (e.g. "35.331.09",
    # "0.0₃0828") already land as text, so no prefix is needed there.
    if ($cellRef -like 'D*' -and $text -match '^\s*-?\d+(\.\d+)?\s*$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}
